# Shift the per-day date header row (B1:AF1) forward by 61 days
# (roster week-view rolled from Aug 2025 -> Oct 2025), then move the
# active selection to I2, which also clears the stale "topLeftCell"
# scroll-position left over from the previous selection at X4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 2; $col -le 32; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Value2 + 61
}

$ws.Range("I2").Select() | Out-Null
